$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 3699.5454
$ws.Cells.Item(98, 9).Value = 1462.125
$ws.Cells.Item(98, 10).Value = 9666
$ws.Cells.Item(98, 11).Value = 1462.125
$ws.Cells.Item(98, 12).Value = 9666
$ws.Cells.Item(98, 13).Value = 35.875
$ws.Cells.Item(98, 14).Value = -12662
$ws.Cells.Item(99, 8).Value = 66678196
$ws.Cells.Item(99, 9).Value = 17006.8
$ws.Cells.Item(99, 10).Value = 200000580
$ws.Cells.Item(99, 11).Value = 51020.39999999999
$ws.Cells.Item(99, 12).Value = 600001740
$ws.Cells.Item(99, 13).Value = -49522.39999999999
$ws.Cells.Item(99, 14).Value = -600004736
$ws.Cells.Item(103, 8).Value = 166668370
$ws.Cells.Item(103, 9).Value = 0
$ws.Cells.Item(103, 11).Value = 0
$ws.Cells.Item(103, 13).ClearContents()
$ws.Cells.Item(106, 8).Value = 10131.167
$ws.Cells.Item(106, 9).Value = 1760.4286
$ws.Cells.Item(106, 11).Value = 1760.4286
$ws.Cells.Item(106, 13).Value = -1129.4286
$ws.Cells.Item(107, 8).Value = 3050.6
$ws.Cells.Item(107, 9).Value = 1277.25
$ws.Cells.Item(107, 11).Value = 1277.25
$ws.Cells.Item(107, 13).Value = 642.75
$ws.Cells.Item(115, 8).Value = 35717868
$ws.Cells.Item(115, 9).Value = 50004516
$ws.Cells.Item(115, 10).Value = 1249.5
$ws.Cells.Item(115, 11).Value = 150013548
$ws.Cells.Item(115, 12).Value = 3748.5
$ws.Cells.Item(115, 13).Value = -150011981
$ws.Cells.Item(115, 14).Value = -6882.5
$ws.Cells.Item(122, 8).Value = 3699.5454
$ws.Cells.Item(122, 9).Value = 1462.125
$ws.Cells.Item(122, 10).Value = 9666
$ws.Cells.Item(122, 11).Value = 4386.375
$ws.Cells.Item(122, 12).Value = 28998
$ws.Cells.Item(122, 13).Value = -1936.375
$ws.Cells.Item(122, 14).Value = -33898
$ws.Cells.Item(127, 8).Value = 909.5
$ws.Cells.Item(127, 10).Value = 909.5
$ws.Cells.Item(127, 12).Value = 2728.5
$ws.Cells.Item(127, 14).Value = -12648.5
$ws.Cells.Item(129, 8).Value = 1922.3182
$ws.Cells.Item(129, 9).Value = 699.63635
$ws.Cells.Item(129, 10).Value = 3145
$ws.Cells.Item(129, 11).Value = 2098.90905
$ws.Cells.Item(129, 12).Value = 9435
$ws.Cells.Item(129, 13).Value = 2901.09095
$ws.Cells.Item(129, 14).Value = -19435
$ws.Cells.Item(135, 8).Value = 1616.4445
$ws.Cells.Item(135, 9).Value = 1796.4
$ws.Cells.Item(135, 10).Value = 716.6667
$ws.Cells.Item(135, 11).Value = 16167.6
$ws.Cells.Item(135, 12).Value = 6450.0003
$ws.Cells.Item(135, 13).Value = -13632.6
$ws.Cells.Item(135, 14).Value = -11520.0003
$ws.Cells.Item(137, 8).Value = 1183.2106
$ws.Cells.Item(137, 9).Value = 1078.8
$ws.Cells.Item(137, 11).Value = 3236.4
$ws.Cells.Item(137, 13).Value = -686.3999999999996
$ws.Cells.Item(138, 8).Value = 2045.4131
$ws.Cells.Item(138, 9).Value = 1266.2
$ws.Cells.Item(138, 10).Value = 2644.8076
$ws.Cells.Item(138, 11).Value = 3798.6
$ws.Cells.Item(138, 12).Value = 7934.4228
$ws.Cells.Item(138, 13).Value = 1341.4
$ws.Cells.Item(138, 14).Value = -18214.4228
$ws.Cells.Item(141, 8).Value = 2994.3
$ws.Cells.Item(141, 9).Value = 3104.7778
$ws.Cells.Item(141, 10).Value = 2000
$ws.Cells.Item(141, 11).Value = 9314.3334
$ws.Cells.Item(141, 12).Value = 6000
$ws.Cells.Item(141, 13).Value = -4134.3334
$ws.Cells.Item(141, 14).Value = -16360
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1888
$ws.Cells.Item(61, 9).Value = 1306.5
$ws.Cells.Item(61, 10).Value = 3283.6
$ws.Cells.Item(61, 11).Value = 1306.5
$ws.Cells.Item(61, 12).Value = 3283.6
$ws.Cells.Item(61, 13).Value = -1094.5
$ws.Cells.Item(61, 14).Value = -3707.6
$ws.Cells.Item(122, 8).Value = 1589.375
$ws.Cells.Item(122, 9).Value = 1532.7142
$ws.Cells.Item(122, 11).Value = 4598.142599999999
$ws.Cells.Item(122, 13).Value = -2148.142599999999
$ws.Cells.Item(136, 8).Value = 1888
$ws.Cells.Item(136, 9).Value = 1306.5
$ws.Cells.Item(136, 10).Value = 3283.6
$ws.Cells.Item(136, 11).Value = 3919.5
$ws.Cells.Item(136, 12).Value = 9850.799999999999
$ws.Cells.Item(136, 13).Value = -1369.5
$ws.Cells.Item(136, 14).Value = -14950.8
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 51610.3
$ws.Cells.Item(107, 9).Value = 59453.41
$ws.Cells.Item(107, 11).Value = 59453.41
$ws.Cells.Item(107, 13).Value = -57533.41
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 722.1667
$ws.Cells.Item(22, 9).Value = 466.6
$ws.Cells.Item(22, 10).Value = 2000
$ws.Cells.Item(22, 11).Value = 466.6
$ws.Cells.Item(22, 12).Value = 2000
$ws.Cells.Item(22, 13).Value = -116.6
$ws.Cells.Item(22, 14).Value = -2700
$ws.Cells.Item(50, 8).Value = 14999.857
$ws.Cells.Item(50, 10).Value = 14999.857
$ws.Cells.Item(50, 12).Value = 14999.857
$ws.Cells.Item(50, 14).Value = -16249.857
$ws.Cells.Item(51, 8).Value = 14999.857
$ws.Cells.Item(51, 10).Value = 14999.857
$ws.Cells.Item(51, 12).Value = 14999.857
$ws.Cells.Item(51, 14).Value = -16471.857
$ws.Cells.Item(59, 8).Value = 19998.75
$ws.Cells.Item(60, 8).Value = 15000
$ws.Cells.Item(60, 10).Value = 15000
$ws.Cells.Item(60, 12).Value = 15000
$ws.Cells.Item(60, 14).Value = -16022
$ws.Cells.Item(61, 8).Value = 14999.857
$ws.Cells.Item(61, 10).Value = 14999.857
$ws.Cells.Item(61, 12).Value = 14999.857
$ws.Cells.Item(61, 14).Value = -15695.857
$ws.Cells.Item(62, 8).Value = 4000
$ws.Cells.Item(62, 9).Value = 4000
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 4000
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).Value = -3376
$ws.Cells.Item(62, 14).ClearContents()
$ws.Cells.Item(65, 8).Value = 4000
$ws.Cells.Item(65, 9).Value = 4000
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 20000
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = -16880
$ws.Cells.Item(65, 14).ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(15, 8).Value = 96.61539
$ws.Cells.Item(15, 9).Value = 47.88889
$ws.Cells.Item(15, 10).Value = 206.25
$ws.Cells.Item(15, 11).Value = 143.66667
$ws.Cells.Item(15, 12).Value = 618.75
$ws.Cells.Item(15, 13).Value = -3.666670000000011
$ws.Cells.Item(15, 14).Value = -898.75
$ws.Cells.Item(18, 8).Value = 92.25
$ws.Cells.Item(18, 9).Value = 92.25
$ws.Cells.Item(18, 11).Value = 276.75
$ws.Cells.Item(18, 13).Value = -107.75
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 14).ClearContents()
$ws.Cells.Item(21, 8).Value = 999
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 14).ClearContents()
$ws.Cells.Item(80, 8).Value = 3777.6667
$ws.Cells.Item(80, 10).Value = 3777.6667
$ws.Cells.Item(80, 12).Value = 11333.0001
$ws.Cells.Item(80, 14).Value = -13205.0001
$ws.Cells.Item(83, 8).Value = 3777.6667
$ws.Cells.Item(83, 10).Value = 3777.6667
$ws.Cells.Item(83, 12).Value = 33999.0003
$ws.Cells.Item(83, 14).Value = -43359.0003
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(74, 8).Value = 39999.332
$ws.Cells.Item(74, 10).Value = 39999.332
$ws.Cells.Item(74, 12).Value = 39999.332
$ws.Cells.Item(74, 14).Value = -41871.332
$ws.Cells.Item(77, 8).Value = 39999.332
$ws.Cells.Item(77, 10).Value = 39999.332
$ws.Cells.Item(77, 12).Value = 119997.996
$ws.Cells.Item(77, 14).Value = -129357.996
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 75002
$ws.Cells.Item(40, 10).Value = 50000
$ws.Cells.Item(40, 12).Value = 50000
$ws.Cells.Item(40, 14).Value = -50272
$ws.Cells.Item(60, 8).Value = 64999.332
$ws.Cells.Item(60, 10).Value = 64999.332
$ws.Cells.Item(60, 12).Value = 64999.332
$ws.Cells.Item(60, 14).Value = -66017.33199999999
$ws.Cells.Item(61, 8).Value = 65252
$ws.Cells.Item(61, 9).Value = 51095.1
$ws.Cells.Item(61, 11).Value = 51095.1
$ws.Cells.Item(61, 13).Value = -50893.1
$ws.Cells.Item(93, 8).Value = 28204.77
$ws.Cells.Item(93, 9).Value = 2777.4167
$ws.Cells.Item(93, 10).Value = 333333
$ws.Cells.Item(93, 11).Value = 2777.4167
$ws.Cells.Item(93, 12).Value = 333333
$ws.Cells.Item(93, 13).Value = -1529.4167
$ws.Cells.Item(93, 14).Value = -335829
$ws.Cells.Item(113, 8).Value = 65252
$ws.Cells.Item(113, 9).Value = 51095.1
$ws.Cells.Item(113, 11).Value = 51095.1
$ws.Cells.Item(113, 13).Value = -48925.1
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1873.8572
$ws.Cells.Item(122, 10).Value = 1673.5
$ws.Cells.Item(122, 12).Value = 5020.5
$ws.Cells.Item(122, 14).Value = -9920.5
$ws.Cells.Item(138, 8).Value = 42000
$ws.Cells.Item(138, 10).Value = 42000
$ws.Cells.Item(138, 12).Value = 42000
$ws.Cells.Item(138, 14).Value = -52280
